$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Ambiente (A3) and URL (B3) to drop the "i-" prefix from the hostname
$ws.Range("A3").Value = "preproducciongestion.segurossura.com.ar"
$ws.Range("B3").Value = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"

# Keep the hyperlink target on B3 in sync with the new URL text
if ($ws.Hyperlinks.Count -gt 0) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address -eq $ws.Range("B3").Address) {
            $hl.Address = "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do"
        }
    }
}

# Move the active selection from C4 to B4
$ws.Range("B4").Select()
